# MIGS.ba.built.4.0.xlsx — add the "culture_collection" field.
#
# The header row (row 15) lists MIGS attribute names alphabetically from
# column H onward. "culture_collection" belongs right after
# "biotic_relationship" (column AG) and before "dew_point" (column AH), i.e.
# it becomes the new column AH. Every field that used to occupy AH..BH
# (dew_point .. trophic_level) therefore shifts one column to the right, into
# AI..BI, taking its header text, cell style and cell comment along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstShiftCol = 34   # column AH (currently "dew_point")
$lastShiftCol  = 60   # column BH (currently "trophic_level")

# Walk from the last column backwards so each destination cell is written
# only after its old contents have already been read.
for ($col = $lastShiftCol; $col -ge $firstShiftCol; $col--) {
    $src = $ws.Cells.Item(15, $col)
    $dst = $ws.Cells.Item(15, $col + 1)

    $srcValue = $src.Value()
    $srcCommentText = $src.Comment.Text()

    $dst.Value = $srcValue

    # Carry the fill/border/font formatting along with the value.
    $src.Copy()
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    if ($col -eq $lastShiftCol) {
        # BI15 is a brand-new cell (the sheet used to stop at BH) so it has
        # no pre-existing comment to edit yet.
        $dst.AddComment($srcCommentText)
    } else {
        # AI15..BH15 already carry the (now one-column-stale) comment that
        # used to sit one column to their left; just refresh its text.
        $dst.Comment.Text($srcCommentText)
    }
}

# Finally, turn what used to be AH15 ("dew_point") into the new
# "culture_collection" header cell, with its own definition comment.
$ah = $ws.Cells.Item(15, $firstShiftCol)
$ah.Value = "culture_collection"
$ah.Comment.Text("Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier")

Write-Host "Inserted culture_collection as new column AH (15 row header); shifted AH..BH -> AI..BI"
